# Update the "엑스포4단지" entry (row 24) to the renamed "엑스포" complex
# with its new code. Excel will automatically rebuild the shared-strings
# table on save: since "엑스포4단지" becomes unreferenced it drops out of
# the table (shifting later indices down by one) while the new "엑스포"
# string is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A24").Value = 184337
$ws.Range("B24").Value = "엑스포"

# Match the author's final selection / view state.
$ws.Range("H20").Select()
